$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.850.23'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '2.898.01'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'" + '366.45'
$ws.Range('E5').Value = '  +4.45%  '
$ws.Range('D6').Value = "'" + '101.22'
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('D7').Value = "'" + '0.539'
$ws.Range('E7').Value = '  -2.87%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -4.22%  '
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').Value = "'" + '18.16'
$ws.Range('E13').Value = '  -3.86%  '
$ws.Range('D14').Value = '3.351.61'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').Value = '2.894.82'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('E17').Value = '  -4.86%  '
$ws.Range('D18').Value = '50.821.65'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').Value = "'" + '3.19'
$ws.Range('E19').Value = '  -5.68%  '
$ws.Range('E20').Value = '  -3.65%  '
$ws.Range('D21').Value = "'" + '12.75'
$ws.Range('E21').Value = '  -4.96%  '
$ws.Range('D22').Value = '0.0₃0936'
$ws.Range('D23').Value = "'" + '67.81'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').Value = "'" + '256.86'
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('E25').Value = '  -1.74%  '
$ws.Range('D26').Value = "'" + '4.32'
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -4.14%  '
$ws.Range('D29').Value = "'" + '25.38'
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').Value = "'" + '6.97'
$ws.Range('E30').Value = '  -5.36%  '
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').Value = "'" + '6.12'
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('E33').Value = '  -3.99%  '
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').Value = "'" + '50.74'
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').Value = "'" + '33.83'
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').Value = "'" + '0.0417'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('E39').Value = '  -5.35%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = "'" + '16.86'
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'" + '2.59'
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('E42').Value = '  -5.99%  '
$ws.Range('E43').Value = '  -3.41%  '
$ws.Range('D44').Value = "'" + '118.96'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = "'" + '21.65'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('D46').Value = "'" + '2.09'
$ws.Range('E46').Value = '  -1.59%  '
$ws.Range('D47').Value = "'" + '2.32'
$ws.Range('E47').Value = '  +1.16%  '
$ws.Range('D48').Value = '2.007.15'
$ws.Range('E48').Value = '  -4.17%  '
$ws.Range('D49').Value = "'" + '3.11'
$ws.Range('E49').Value = '  -5.82%  '
$ws.Range('D50').Value = '3.181.78'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('E51').Value = '  -1.97%  '
